# Excel COM-interop script: add "2022-Q4" sheet + update "总计" (totals) sheet.
#
# Target layout after edit:
#   Sheets: 总计, 2022-Q4 (new), 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4
#   总计: a new row 2 is inserted for 2022-Q4 (14 holdings, 2.32 亿元),
#         existing rows shift down by one.
#   2022-Q4: brand-new worksheet with the same layout as the other
#            quarterly fund-holding sheets (basic fund info table).
#
# Note: this runtime's COM bridge does not deliver Copy()/PasteSpecial()
# content into a worksheet created in the same script run via
# Worksheets.Add() (format/values silently fail to land), so the new
# quarter sheet is produced by *duplicating* an existing quarterly sheet
# (Worksheet.Copy, which faithfully clones data + styles) and then
# overwriting its contents - that keeps every style index identical to
# the ones already used elsewhere in the workbook instead of minting new
# near-duplicate styles.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (totals) sheet: push existing data rows down by
#    one row and write the new 2022-Q4 summary row at the top.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item(1)

# Grow the formatting of row 6 (new last row) from row 5 first, so the
# freshly-used row inherits the same per-column look (bold/bordered A
# column, plain B:D columns) as every other data row already has.
$totals.Range("A5:D5").Copy()
$totals.Range("A6:D6").PasteSpecial(-4122)  # xlPasteFormats

# Shift the four existing quarters down one row (bottom-up so we never
# clobber a source row before it has been copied down). Value2 (not
# Value) is required when chaining a read straight into another cell's
# setter here. Column A is not "data" - it is always just the 0-based
# row index, so it is (re)written with a literal number rather than
# carried over from the source row.
$totals.Cells.Item(6,2).Value = $totals.Cells.Item(5,2).Value2
$totals.Cells.Item(6,3).Value = $totals.Cells.Item(5,3).Value2
$totals.Cells.Item(6,4).Value = $totals.Cells.Item(5,4).Value2
$totals.Cells.Item(6,1).Value = 4

$totals.Cells.Item(5,2).Value = $totals.Cells.Item(4,2).Value2
$totals.Cells.Item(5,3).Value = $totals.Cells.Item(4,3).Value2
$totals.Cells.Item(5,4).Value = $totals.Cells.Item(4,4).Value2
$totals.Cells.Item(5,1).Value = 3

$totals.Cells.Item(4,2).Value = $totals.Cells.Item(3,2).Value2
$totals.Cells.Item(4,3).Value = $totals.Cells.Item(3,3).Value2
$totals.Cells.Item(4,4).Value = $totals.Cells.Item(3,4).Value2
$totals.Cells.Item(4,1).Value = 2

$totals.Cells.Item(3,2).Value = $totals.Cells.Item(2,2).Value2
$totals.Cells.Item(3,3).Value = $totals.Cells.Item(2,3).Value2
$totals.Cells.Item(3,4).Value = $totals.Cells.Item(2,4).Value2
$totals.Cells.Item(3,1).Value = 1

# New top data row: 2022-Q4.
$totals.Cells.Item(2,1).Value = 0
$totals.Cells.Item(2,2).Value = "2022-Q4"
$totals.Cells.Item(2,3).Value = 14
$totals.Cells.Item(2,4).Value = 2.32

# ---------------------------------------------------------------------
# 2) Insert the brand-new "2022-Q4" worksheet right after "总计" by
#    duplicating "2022-Q3" (identical table layout/styles) and then
#    overwriting its contents in place.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q3")
$template.Copy($null, $totals)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The template only has 7 data rows (1 header + 6 funds); the new sheet
# needs 15 (1 header + 14 funds), so clone row 7's per-column formatting
# down through row 15 before writing the extra rows' values.
$q4.Range("A7:H7").Copy()
$q4.Range("A8:H15").PasteSpecial(-4122)  # xlPasteFormats

# -- header row (unchanged text, rewritten for clarity/self-containment) -
$q4.Cells.Item(1,2).Value = "基金代码"
$q4.Cells.Item(1,3).Value = "基金名称"
$q4.Cells.Item(1,4).Value = "基金规模"
$q4.Cells.Item(1,5).Value = "股票总仓位"
$q4.Cells.Item(1,6).Value = "仓位占比"
$q4.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4.Cells.Item(1,8).Value = "仓位排名"

# -- data rows -----------------------------------------------------------
# Columns: A idx(n) | B code(text) | C name(text) | D scale(text) |
#          E position(text) | F ratio(text) | G value(text, mostly) |
#          H rank(n)
# B/D/E/F/G hold digit-only strings (fund codes with leading zeros,
# percentages, etc.) that Excel would otherwise silently coerce to
# numbers, so they are written as `="literal"` formulas first and then
# flattened to static values with a copy/paste-special (values only).
# That keeps them genuine text cells without adding any new style.

$rows = @(
    @(0,  "920003", "中金新锐股票A",                     "19.30", "91.44", "6.90", "1.3317", 2),
    @(1,  "860001", "光大阳光混合A",                     "9.17",  "87.55", "4.66", "0.4273", 5),
    @(2,  "920923", "中金新锐股票C",                     "3.33",  "91.44", "6.90", "0.2298", 2),
    @(3,  "920002", "中金精选股票A",                     "3.21",  "93.08", "3.80", "0.1220", 3),
    @(4,  "860052", "光大阳光启明星创新驱动主题混合B",   "1.73",  "87.94", "5.20", "0.0900", 3),
    @(5,  "860053", "光大阳光启明星创新驱动主题混合C",   "1.26",  "87.94", "5.20", "0.0655", 3),
    @(6,  "860016", "光大阳光启明星创新驱动主题混合A",   "0.45",  "87.94", "5.20", "0.0234", 3),
    @(7,  "002236", "大成中证360互联网+大数据100指数A", "1.15",  "92.50", "0.99", "0.0114", 9),
    @(8,  "003359", "大成中证360互联网+大数据100指数C", "1.12",  "92.50", "0.99", "0.0111", 9),
    @(9,  "860036", "光大阳光混合B",                     "0.10",  "87.55", "4.66", "0.0047", 5),
    @(10, "920922", "中金精选股票C",                     "0.12",  "93.08", "3.80", "0.0046", 3),
    @(11, "005095", "国泰量化成长优选混合A",             "0.20",  "87.25", "1.49", "0.0030", 10),
    @(12, "005096", "国泰量化成长优选混合C",             "0.02",  "87.25", "1.49", "0.0003", 10),
    @(13, "860037", "光大阳光混合C",                     "0.00",  "87.55", "4.66", $null,     5)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $q4.Cells.Item($r,1).Value = $row[0]
    $q4.Cells.Item($r,2).Formula = '="' + $row[1] + '"'
    $q4.Cells.Item($r,3).Value = $row[2]
    $q4.Cells.Item($r,4).Formula = '="' + $row[3] + '"'
    $q4.Cells.Item($r,5).Formula = '="' + $row[4] + '"'
    $q4.Cells.Item($r,6).Formula = '="' + $row[5] + '"'
    if ($row[6] -ne $null) {
        $q4.Cells.Item($r,7).Formula = '="' + $row[6] + '"'
    } else {
        $q4.Cells.Item($r,7).Value = 0
    }
    $q4.Cells.Item($r,8).Value = $row[7]
}

# Flatten the helper `="text"` formulas above into plain static text
# values (covering every row written, column by column).
$textRange = $q4.Range("B2:B15")
$textRange.Copy()
$textRange.PasteSpecial(-4163)  # xlPasteValues

$textRange = $q4.Range("D2:G14")
$textRange.Copy()
$textRange.PasteSpecial(-4163)  # xlPasteValues

$textRange = $q4.Range("D15:F15")
$textRange.Copy()
$textRange.PasteSpecial(-4163)  # xlPasteValues
